# Locate the trailing empty "ListParagraph" paragraph at the very end of
# the body (immediately before sectPr) -- this is the insertion point the
# new "What is an Exception?" / SQL Exceptions / SQL Injection Attacks /
# Cons of the Statement Object section is built from.
$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$r = $p.Range

# Sanity-check we grabbed the right (empty, ListParagraph-styled) paragraph
# before clobbering it. Range.Text includes the trailing paragraph mark
# (CR, chr 13), so strip that before checking for "empty".
$bodyText = $r.Text.TrimEnd([char]13)
if ($bodyText -ne "" -or $p.Range.ParagraphStyle.NameLocal -ne "List Paragraph") {
    throw "Unexpected final paragraph -- refusing to overwrite."
}

# Replacing this paragraph's Range with a run of <w:p> siblings swaps the
# single empty paragraph out for the whole new block in one shot, while
# keeping it anchored right before </w:body>/<w:sectPr> the same way the
# original paragraph was. The first <w:p> keeps the original paragraph's
# identity (w14:paraId=515E3CE7 etc.) since the diff shows that opening tag
# as unchanged context -- only its <w:pPr>/content changed.
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="515E3CE7" w14:textId="77777777" w:rsidR="00E5432C" w:rsidRPr="001B2218" w:rsidRDefault="00E5432C" w:rsidP="001B2218"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>What is an Exception?</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>An Exception is an Error that crashes your program.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">An Exception can be thrown using a </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>throws</w:t></w:r><w:r><w:t xml:space="preserve"> clause, or caught using a </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>try / catch block</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Some programming statements require an Exception to be thrown or caught.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Throwing an Exception does not prevent your program from </w:t></w:r><w:r><w:t>crashing but</w:t></w:r><w:r><w:t xml:space="preserve"> catching an Exception does. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>SQL</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> Exception</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>s</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">An </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>SQLException</w:t></w:r><w:r><w:t xml:space="preserve"> occurs if a SQL statement is written incorrectly.</w:t></w:r></w:p><w:p/><w:p/><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">SQL </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Injection Attacks</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">A SQL Injection Attack occurs when SWL statements entered through an application to expose or harm a database. </w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Companies that have been affected by SQL injection attacks are Sony, TJX, MasterCard, PBS, etc. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">SQL Injection Attacks can occur by incorrectly handling escape characters such as quotes, assigning conditions as input, etc. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Cons of the Statement Object</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">The </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Statement</w:t></w:r><w:r><w:t xml:space="preserve"> Object does not escape quotes.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Writing </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>INSERT</w:t></w:r><w:r><w:t xml:space="preserve"> statements requires </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>concatenation</w:t></w:r><w:r><w:t xml:space="preserve"> making it tedious to put one together. </w:t></w:r></w:p>
'@

$r.InsertXML($xml)

Write-Output "Inserted Exceptions / SQL Exceptions / SQL Injection Attacks / Cons of the Statement Object section."
